$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDiBCpDoC")
$ws.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"
